# Update the "Förändrad" (Changed) date column (C) for rows 2-7
# from serial date 45184 (2023-09-15) to 45185 (2023-09-16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = 45185
}
